# Dataset Alkohol.xlsx - remove trailing " A"/" B" placeholder suffixes from
# product_name values (column E) in Sheet1, and leave the selection on the
# last edited cell (E23), matching the author's cleanup commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value  = "Red Wine "
$ws.Range("E7").Value  = "White Wine"
$ws.Range("E8").Value  = "Whiskey "
$ws.Range("E9").Value  = "Vodka "
$ws.Range("E10").Value = "Rum "
$ws.Range("E11").Value = "Gin "
$ws.Range("E12").Value = "Cider "
$ws.Range("E13").Value = "Sake "
$ws.Range("E14").Value = "Cocktail RTD "
$ws.Range("E15").Value = "Cocktail RTD "
$ws.Range("E16").Value = "Sparkling Wine "
$ws.Range("E17").Value = "Craft Beer IPA "
$ws.Range("E18").Value = "Craft Beer Stout "
$ws.Range("E19").Value = "Tequila "
$ws.Range("E20").Value = "Mead "
$ws.Range("E21").Value = "Rice Wine "
$ws.Range("E22").Value = "Palm Wine "
$ws.Range("E23").Value = "Low Alcohol Beer "

$ws.Range("E23").Select()
